{"js": "// The \"Form of consent for participants unable to provide a signature or to\n// mark the box:\" paragraph ends with two runs: an italic black space (\" \")\n// followed by an italic red placeholder (\"N/A or details here\"). The edit\n// collapses both into a single, non-italic, black run reading \" N/A\".\nconst body = context.document.body;\n\n// Search across the run boundary so the hit covers both the trailing space\n// run and the placeholder run that follows it.\nconst results = body.search(\" N/A or details here\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the placeholder text to update.\");\n}\n\nconst target = results.items[0];\n\n// Normalize formatting to match the surrounding non-italic black text...\ntarget.font.italic = false;\ntarget.font.color = \"#000000\";\n\n// ...and replace the whole matched range's text with the shortened \" N/A\".\ntarget.insertText(\" N/A\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The \"Form of consent for participants unable to provide a signature or to\n# mark the box:\" paragraph ends with two runs: an italic black space (\" \")\n# followed by an italic red placeholder (\"N/A or details here\"). The edit\n# collapses both into a single, non-italic, black run reading \" N/A\".\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n# MatchCase so we land on the exact placeholder span (covers both the\n# trailing-space run and the \"N/A or details here\" run that follows it).\n$found = $rng.Find.Execute(\" N/A or details here\", $true)\n\nif (-not $found) {\n    throw \"Could not find the placeholder text to update.\"\n}\n\n# Normalize formatting to match the surrounding non-italic black text...\n$rng.Font.Italic = $false\n$rng.Font.Color = 0\n\n# ...and replace the whole matched range's text with the shortened \" N/A\".\n$rng.Text = \" N/A\"\n"}
